$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Common")
Write-Host $ws.Name
Write-Host $ws.Range("A70").Value
